$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update B5 value from "category_one" to "similar_category_one"
$ws.Range("B5").Value = "similar_category_one"

# Delete column D entirely (removes the "similar_match" header column and its data)
$ws.Range("D1:D5").Delete()

# Update the selected cell/range shown in the sheet view
$ws.Range("K7").Select()
